$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53-180 down to 54-181
$ws.Rows("53:53").Insert()

# Populate the new row 53 with the new weekly record
$ws.Range("A53").Value = 4
$ws.Range("B53").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C53").Value = "Los Lagos"
$ws.Range("D53").Value = 44526
$ws.Range("E53").Value = 10
$ws.Range("F53").Value = 100112003
$ws.Range("G53").Value = "Ajo"
$ws.Range("H53").Value = "Chino"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 240
$ws.Range("K53").Value = 20000
$ws.Range("L53").Value = 21000
$ws.Range("M53").Value = 20500
$ws.Range("N53").Value = "$/caja 10 kilos"
$ws.Range("O53").Value = "China"
$ws.Range("P53").Value = 2050
$ws.Range("Q53").Value = 10
$ws.Range("R53").Value = "Hortaliza"
